# Refresh the coinranking.com cryptocurrency price/volume snapshot.
# (GitHub Actions scheduled update - Thu Oct 26 15:38:40 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.004.19"
$ws.Range("E2").Value = "  -2.30%  "
$ws.Range("D3").Value = "1.776.06"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.52"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.71%  "
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.25"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.285"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0703"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.71%  "
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("D12").Value = "2.034.24"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "1.788.06"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.54"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.98%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.620"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.95%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "33.940.45"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.19"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.61"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.65"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -6.04%  "
$ws.Range("D20").Value = "0.0₃0768"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.49"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.04"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.18%  "
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.58"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.25"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.97"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("E28").Value = "  -3.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0515"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.66"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.40%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.19"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.50"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.92%  "
$ws.Range("D35").Value = "1.392.33"
$ws.Range("E35").Value = "  -4.74%  "
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.625"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0185"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.927"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "78.80"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -6.02%  "
$ws.Range("E42").Value = "  -5.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.09"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.82"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.93%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.04"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("D47").Value = "1.926.31"
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.40"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.995"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.72"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.67%  "
$ws.Range("D51").Value = "0.0₆0117"
$ws.Range("E51").Value = "  -3.79%  "
